$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.3254
$ws.Cells.Item(2, 3).Value = 0.3254
$ws.Cells.Item(2, 4).Value = 0.3086000084877014
$ws.Cells.Item(2, 5).Value = 0.5224000215530396
$ws.Cells.Item(2, 6).Value = 0.1756999939680099

$ws.Cells.Item(3, 2).Value = 0.3977
$ws.Cells.Item(3, 3).Value = 0.3977
$ws.Cells.Item(3, 4).Value = 0.4076
$ws.Cells.Item(3, 5).Value = 0.6488999724388123
$ws.Cells.Item(3, 6).Value = 0.1066000014543533

$ws.Cells.Item(4, 2).Value = 0.3857
$ws.Cells.Item(4, 3).Value = 0.3857
$ws.Cells.Item(4, 4).Value = 0.3877
$ws.Cells.Item(4, 5).Value = 0.503600001335144
$ws.Cells.Item(4, 6).Value = 0.1861999928951263

$ws.Cells.Item(5, 2).Value = 0.2531
$ws.Cells.Item(5, 3).Value = 0.2531
$ws.Cells.Item(5, 4).Value = 0.2237
$ws.Cells.Item(5, 5).Value = 0.4767000079154968
$ws.Cells.Item(5, 6).Value = 0.09839999675750732

$ws.Cells.Item(6, 2).Value = 0.2696
$ws.Cells.Item(6, 3).Value = 0.2696
$ws.Cells.Item(6, 4).Value = 0.2793
$ws.Cells.Item(6, 5).Value = 0.4860999882221222
$ws.Cells.Item(6, 6).Value = 0.1142000034451485

$ws.Cells.Item(7, 2).Value = 0.2953
$ws.Cells.Item(7, 3).Value = 0.2953
$ws.Cells.Item(7, 4).Value = 0.3003999888896942
$ws.Cells.Item(7, 5).Value = 0.3982000052928925
$ws.Cells.Item(7, 6).Value = 0.1956000030040741

$ws.Cells.Item(8, 2).Value = 0.2883
$ws.Cells.Item(8, 3).Value = 0.2883
$ws.Cells.Item(8, 4).Value = 0.287
$ws.Cells.Item(8, 5).Value = 0.4081999957561493
$ws.Cells.Item(8, 6).Value = 0.1956000030040741
